# Insert a new data row at row 427 (a new weekly price observation for
# Femacal de La Calera - Ajo). Excel's Rows.Insert() pushes the existing
# row 427 (and everything below it) down by one, copying row formatting
# from the row above - so we follow up by copying the now-shifted row 428
# (the former row 427) back into the new row 427 to duplicate all of its
# values, then overwrite just the six cells (Fecha, Volumen, Precio
# minimo/maximo/promedio ponderado, Precio $/Kg) that differ for the new
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new row; everything from row 427 down shifts to 428..455
$ws.Rows.Item(427).Insert()

# Duplicate the (now shifted) original row 427 data - currently sitting in
# row 428 - into the freshly inserted, blank row 427.
$ws.Range("A428:R428").Copy()
$ws.Range("A427").PasteSpecial()

# Overwrite the new row 427 with this week's actual figures.
$ws.Range("D427").Value = 44714
$ws.Range("J427").Value = 76
$ws.Range("K427").Value = 17000
$ws.Range("L427").Value = 18000
$ws.Range("M427").Value = 17500
$ws.Range("P427").Value = 1750
